$wb = $excel.ActiveWorkbook

# DATASHEET: update values and move selection to F3
$wsDataSheet = $wb.Worksheets.Item("DATASHEET")
$wsDataSheet.Range("E3").Value = 20
$wsDataSheet.Range("F3").Value = 40
$wsDataSheet.Activate()
$wsDataSheet.Range("F3").Select()

# MOBILE_CONFIGURATION: move selection to E3
$wsMobileConfig = $wb.Worksheets.Item("MOBILE_CONFIGURATION")
$wsMobileConfig.Activate()
$wsMobileConfig.Range("E3").Select()

# MAIL_SEND: change B2 value from "Y" to "N" and move selection to B2
$wsMailSend = $wb.Worksheets.Item("MAIL_SEND")
$wsMailSend.Range("B2").Value = "N"
$wsMailSend.Activate()
$wsMailSend.Range("B2").Select()

# MAIN_CONTROLLER: becomes the active/selected tab, selection stays at B3
$wsMainController = $wb.Worksheets.Item("MAIN_CONTROLLER")
$wsMainController.Activate()
$wsMainController.Range("B3").Select()
